$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8302551507949829
$ws.Range("B1").Value = 2.177186727523804
$ws.Range("C1").Value = 4.904501914978027
$ws.Range("D1").Value = 2.391692638397217
$ws.Range("E1").Value = 1.356216073036194
